$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two blank rows at row 14 (pushes existing rows 14-121 down to 16-123)
$ws.Rows("14:15").Insert()

# 2. Fill in the new bottom rows (appended after the old last row, now row 123)
#    Order matters for shared-string table layout, matching the original authoring order.
$ws.Range("A125").Value = "digital_investigation_req_message"
$ws.Range("B125").Value = "You need to flag at least three items from the investigation to proceed."

$ws.Range("A124").Value = "digital_investigation_report"
$ws.Range("B124").Value = "Digital Investigation Report"

$ws.Range("A126").Value = "report"
$ws.Range("B126").Value = "Report"

# 3. Fill the newly inserted rows near the top
$ws.Range("A14").Value = "proceed"
$ws.Range("B14").Value = "PROCEED"

$ws.Range("A15").Value = "cancel"
$ws.Range("B15").Value = "CANCEL"

# 4. Fill the final new row at the bottom
$ws.Range("A127").Value = "digital_investigation_report_confirm"
$ws.Range("B127").Value = "Do you want to finish the report and proceed?"

# 5. Update the view/selection to match the authored state
$ws.Range("B127").Select()
$app.ActiveWindow.ScrollRow = 112
